$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 3886.879
$ws.Cells.Item(17, 10).Value = 4198.9
$ws.Cells.Item(17, 12).Value = 12596.7
$ws.Cells.Item(17, 14).Value = -12932.7
$ws.Cells.Item(28, 8).Value = 51268.5
$ws.Cells.Item(28, 9).Value = 72125.36
$ws.Cells.Item(28, 10).Value = 2602.5
$ws.Cells.Item(28, 11).Value = 72125.36
$ws.Cells.Item(28, 12).Value = 2602.5
$ws.Cells.Item(28, 13).Value = -71640.36
$ws.Cells.Item(28, 14).Value = -3572.5
$ws.Cells.Item(86, 8).Value = 5304.9165
$ws.Cells.Item(86, 10).Value = 6133
$ws.Cells.Item(86, 12).Value = 6133
$ws.Cells.Item(86, 14).Value = -8379
$ws.Cells.Item(89, 8).Value = 5304.9165
$ws.Cells.Item(89, 10).Value = 6133
$ws.Cells.Item(89, 12).Value = 30665
$ws.Cells.Item(89, 14).Value = -41897
$ws.Cells.Item(107, 8).Value = 4024.1
$ws.Cells.Item(107, 9).Value = 4669.615
$ws.Cells.Item(107, 10).Value = 2825.2856
$ws.Cells.Item(107, 11).Value = 4669.615
$ws.Cells.Item(107, 12).Value = 2825.2856
$ws.Cells.Item(107, 13).Value = -2749.615
$ws.Cells.Item(107, 14).Value = -6665.2856

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 27780126
$ws.Cells.Item(74, 9).Value = 33335804
$ws.Cells.Item(74, 10).Value = 1729.5
$ws.Cells.Item(74, 11).Value = 33335804
$ws.Cells.Item(74, 12).Value = 1729.5
$ws.Cells.Item(74, 13).Value = -33334930
$ws.Cells.Item(74, 14).Value = -3477.5
$ws.Cells.Item(77, 8).Value = 27780126
$ws.Cells.Item(77, 9).Value = 33335804
$ws.Cells.Item(77, 10).Value = 1729.5
$ws.Cells.Item(77, 11).Value = 166679020
$ws.Cells.Item(77, 12).Value = 8647.5
$ws.Cells.Item(77, 13).Value = -166674652
$ws.Cells.Item(77, 14).Value = -17383.5
$ws.Cells.Item(97, 8).Value = 956.25
$ws.Cells.Item(97, 9).Value = 1411.3636
$ws.Cells.Item(97, 11).Value = 1411.3636
$ws.Cells.Item(97, 13).Value = -915.3635999999999
$ws.Cells.Item(110, 8).Value = 4088.4443
$ws.Cells.Item(110, 9).Value = 2721
$ws.Cells.Item(110, 10).Value = 5182.4
$ws.Cells.Item(110, 11).Value = 2721
$ws.Cells.Item(110, 12).Value = 5182.4
$ws.Cells.Item(110, 13).Value = -676
$ws.Cells.Item(110, 14).Value = -9272.4
$ws.Cells.Item(122, 8).Value = 1942.6316
$ws.Cells.Item(122, 9).Value = 1339
$ws.Cells.Item(122, 11).Value = 4017
$ws.Cells.Item(122, 13).Value = -1567

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 4709.385
$ws.Cells.Item(20, 9).Value = 3299.2
$ws.Cells.Item(20, 11).Value = 3299.2
$ws.Cells.Item(20, 13).Value = -3052.2

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3809.3
$ws.Cells.Item(16, 9).Value = 3365.5
$ws.Cells.Item(16, 10).Value = 4475
$ws.Cells.Item(16, 11).Value = 3365.5
$ws.Cells.Item(16, 12).Value = 4475
$ws.Cells.Item(16, 13).Value = -3078.5
$ws.Cells.Item(16, 14).Value = -5049
$ws.Cells.Item(62, 8).Value = 5980.1113
$ws.Cells.Item(62, 9).Value = 2606.3333
$ws.Cells.Item(62, 11).Value = 2606.3333
$ws.Cells.Item(62, 13).Value = -1982.3333
$ws.Cells.Item(65, 8).Value = 5980.1113
$ws.Cells.Item(65, 9).Value = 2606.3333
$ws.Cells.Item(65, 11).Value = 13031.6665
$ws.Cells.Item(65, 13).Value = -9911.666499999999
$ws.Cells.Item(113, 8).Value = 3809.3
$ws.Cells.Item(113, 9).Value = 3365.5
$ws.Cells.Item(113, 10).Value = 4475
$ws.Cells.Item(113, 11).Value = 3365.5
$ws.Cells.Item(113, 12).Value = 4475
$ws.Cells.Item(113, 13).Value = -1195.5
$ws.Cells.Item(113, 14).Value = -8815
$ws.Cells.Item(122, 8).Value = 5820.65
$ws.Cells.Item(122, 9).Value = 1672.4286
$ws.Cells.Item(122, 11).Value = 5017.2858
$ws.Cells.Item(122, 13).Value = -2567.2858
$ws.Cells.Item(132, 8).Value = 2294.6843
$ws.Cells.Item(132, 9).Value = 1365.8334
$ws.Cells.Item(132, 11).Value = 4097.5002
$ws.Cells.Item(132, 13).Value = -1567.5002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(52, 8).Value = 376.66666
$ws.Cells.Item(52, 10).Value = 376.66666
$ws.Cells.Item(52, 12).Value = 1129.99998
$ws.Cells.Item(52, 14).Value = -1661.99998
$ws.Cells.Item(68, 8).Value = 1259.5
$ws.Cells.Item(68, 9).Value = 1324.5
$ws.Cells.Item(68, 10).Value = 1216.1666
$ws.Cells.Item(68, 11).Value = 3973.5
$ws.Cells.Item(68, 12).Value = 3648.4998
$ws.Cells.Item(68, 13).Value = -3162.5
$ws.Cells.Item(68, 14).Value = -5270.4998
$ws.Cells.Item(71, 8).Value = 1259.5
$ws.Cells.Item(71, 9).Value = 1324.5
$ws.Cells.Item(71, 10).Value = 1216.1666
$ws.Cells.Item(71, 11).Value = 11920.5
$ws.Cells.Item(71, 12).Value = 10945.4994
$ws.Cells.Item(71, 13).Value = -7864.5
$ws.Cells.Item(71, 14).Value = -19057.4994
$ws.Cells.Item(117, 8).Value = 4139.5
$ws.Cells.Item(117, 9).Value = 508.66666
$ws.Cells.Item(117, 11).Value = 1525.99998
$ws.Cells.Item(117, 13).Value = 1916.00002
$ws.Cells.Item(129, 8).Value = 5955458.5
$ws.Cells.Item(129, 9).Value = 644.375
$ws.Cells.Item(129, 11).Value = 1933.125
$ws.Cells.Item(129, 13).Value = 3066.875
$ws.Cells.Item(132, 8).Value = 2739.2068
$ws.Cells.Item(132, 9).Value = 2048.4614
$ws.Cells.Item(132, 11).Value = 18436.1526
$ws.Cells.Item(132, 13).Value = -15906.1526

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(52, 8).Value = 29999
$ws.Cells.Item(52, 10).Value = 29999
$ws.Cells.Item(52, 12).Value = 29999
$ws.Cells.Item(52, 14).Value = -30517
$ws.Cells.Item(102, 8).Value = 2796.9333
$ws.Cells.Item(102, 9).Value = 2001.2727
$ws.Cells.Item(102, 11).Value = 2001.2727
$ws.Cells.Item(102, 13).Value = -379.2727
$ws.Cells.Item(113, 8).Value = 5043.3687
$ws.Cells.Item(113, 9).Value = 4633.615
$ws.Cells.Item(113, 11).Value = 4633.615
$ws.Cells.Item(113, 13).Value = -2463.615
$ws.Cells.Item(132, 8).Value = 9732.360000000001
$ws.Cells.Item(132, 9).Value = 8286.666999999999
$ws.Cells.Item(132, 10).Value = 11900.9
$ws.Cells.Item(132, 11).Value = 24860.001
$ws.Cells.Item(132, 12).Value = 35702.7
$ws.Cells.Item(132, 13).Value = -22330.001
$ws.Cells.Item(132, 14).Value = -40762.7

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6540.613
$ws.Cells.Item(7, 9).Value = 3758.4119
$ws.Cells.Item(7, 10).Value = 9919
$ws.Cells.Item(7, 11).Value = 3758.4119
$ws.Cells.Item(7, 12).Value = 9919
$ws.Cells.Item(7, 13).Value = -3646.4119
$ws.Cells.Item(7, 14).Value = -10143
$ws.Cells.Item(61, 8).Value = 3583.4285
$ws.Cells.Item(61, 9).Value = 930.25
$ws.Cells.Item(61, 11).Value = 930.25
$ws.Cells.Item(61, 13).Value = -728.25
$ws.Cells.Item(113, 8).Value = 3583.4285
$ws.Cells.Item(113, 9).Value = 930.25
$ws.Cells.Item(113, 11).Value = 930.25
$ws.Cells.Item(113, 13).Value = 1239.75
$ws.Cells.Item(126, 8).Value = 6540.613
$ws.Cells.Item(126, 9).Value = 3758.4119
$ws.Cells.Item(126, 10).Value = 9919
$ws.Cells.Item(126, 11).Value = 11275.2357
$ws.Cells.Item(126, 12).Value = 29757
$ws.Cells.Item(126, 13).Value = -8805.235700000001
$ws.Cells.Item(126, 14).Value = -34697
$ws.Cells.Item(132, 8).Value = 7584.1665
$ws.Cells.Item(132, 9).Value = 4333.3335
$ws.Cells.Item(132, 10).Value = 10835
$ws.Cells.Item(132, 11).Value = 13000.0005
$ws.Cells.Item(132, 12).Value = 32505
$ws.Cells.Item(132, 13).Value = -10470.0005
$ws.Cells.Item(132, 14).Value = -37565
$ws.Cells.Item(136, 8).Value = 11864.643
$ws.Cells.Item(136, 9).Value = 3766.3333
$ws.Cells.Item(136, 10).Value = 17938.375
$ws.Cells.Item(136, 11).Value = 11298.9999
$ws.Cells.Item(136, 12).Value = 53815.125
$ws.Cells.Item(136, 13).Value = -8748.999899999999
$ws.Cells.Item(136, 14).Value = -58915.125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(39, 8).Value = 516524
$ws.Cells.Item(39, 10).Value = 33049
$ws.Cells.Item(39, 12).Value = 33049
$ws.Cells.Item(39, 14).Value = -33875
$ws.Cells.Item(69, 8).Value = 17500
$ws.Cells.Item(69, 10).Value = 17500
$ws.Cells.Item(69, 12).Value = 17500
$ws.Cells.Item(69, 14).Value = -18998
$ws.Cells.Item(72, 8).Value = 17500
$ws.Cells.Item(72, 10).Value = 17500
$ws.Cells.Item(72, 12).Value = 52500
$ws.Cells.Item(72, 14).Value = -59988
$ws.Cells.Item(132, 8).Value = 7887.2974
$ws.Cells.Item(132, 9).Value = 7318.2354
$ws.Cells.Item(132, 10).Value = 14336.667
$ws.Cells.Item(132, 11).Value = 21954.7062
$ws.Cells.Item(132, 12).Value = 43010.001
$ws.Cells.Item(132, 13).Value = -19424.7062
$ws.Cells.Item(132, 14).Value = -48070.001
$ws.Cells.Item(136, 8).Value = 11800.8
$ws.Cells.Item(136, 9).Value = 8333
$ws.Cells.Item(136, 11).Value = 24999
$ws.Cells.Item(136, 13).Value = -22449
